$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test_number column (A2:A11) from 2 to 3
$ws.Range("A2:A11").Value = 3

# Update the active selection to I9
$ws.Range("I9").Select()
